$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "First player"
$ws.Range("C2").Value = 80.0

$ws.Range("A3").Value = 2.0
$ws.Range("B3").Value = "Second player"
$ws.Range("C3").Value = 71.0
